$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Text updates for the "Experiment: ..." section headers
#    (seed labels switch from "seed=N" to "seed_N"; hammett -> Hammett)
# -----------------------------------------------------------------
$ws.Range("A4").Value  = "Experiment: XGBC_seed_1"
$ws.Range("A6").Value  = "Experiment: XGBC_seed_2"
$ws.Range("A8").Value  = "Experiment: XGBC_seed_3"
$ws.Range("A10").Value = "Experiment: XGBC_seed_4"
$ws.Range("A12").Value = "Experiment: XGBC_seed_5"
$ws.Range("A14").Value = "Experiment: Hammett"

# -----------------------------------------------------------------
# 2. Updated numbers in the "Test accuracy (seed = 1)" table
#    (rows 20-24, columns B..I)
# -----------------------------------------------------------------

# Row 21 - RFC
$ws.Range("E21").Value = 0.59
$ws.Range("F21").Value = 0.94

# Row 22 - LogR
$ws.Range("B22").Value = 0.88
$ws.Range("E22").Value = 0.09
$ws.Range("F22").Value = 0.97
$ws.Range("I22").Value = 0.09

# Row 23 - NNC
$ws.Range("B23").Value = 0.78
$ws.Range("D23").Value = 0.91
$ws.Range("E23").Value = 0.03
$ws.Range("F23").Value = 0.84
$ws.Range("G23").Value = 0.88
$ws.Range("I23").Value = 0.06

# Row 24 - GPC
$ws.Range("B24").Value = 0.06
$ws.Range("F24").Value = 0.06

# Bold-highlight changes that accompany the value updates above
# (cells that stop being the row's "best" value lose bold, the new
# best value gains it)
$ws.Range("F20").Font.Bold = $false
$ws.Range("D23").Font.Bold = $false
$ws.Range("G23").Font.Bold = $false
$ws.Range("F22").Font.Bold = $true

# -----------------------------------------------------------------
# 3. Updated numbers in the "Train accuracy (seed = 1)" table
#    (rows 29-33, columns B..I)
# -----------------------------------------------------------------
$ws.Range("B31").Value = 1
$ws.Range("F31").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("F32").Value = 1
$ws.Range("B33").Value = 1
$ws.Range("F33").Value = 1

# New (empty) formatted cell that now exists right below the table
$ws.Range("F34").Value = $null
$ws.Range("F34").NumberFormat = "0.00"
$ws.Range("F34").Font.Bold = $true
$ws.Range("F34").HorizontalAlignment = -4108

# -----------------------------------------------------------------
# 4. Selection / scroll position bookkeeping
# -----------------------------------------------------------------
$ws.Range("D26").Select()
